$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column to the left (M) so the freshly
# inserted column can be given a matching custom width.
$matchWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a new blank column before column N (14th column), shifting
# existing N, O, P columns to O, P, Q respectively.
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $matchWidth

# Update the selection on this sheet to T7
$ws.Range("T7").Select()

# Make the "Repayment schedule" sheet the active one (this becomes the
# selected tab when the workbook is saved)
$ws.Activate()
